# Auto-generated edit script: update crypto price/volume table
# Uses a helper that forces numeric-looking strings to be stored as TEXT
# (mirrors the source data which stores inline strings, not numbers),
# by evaluating TEXT("value","@") as a formula and then pasting the
# formula result back as a literal value (no residual number formatting).
function Set-TextValue($ws, $Address, $Text) {
    $escaped = $Text.Replace('"', '""')
    $range = $ws.Range($Address)
    $range.Formula = '=TEXT("' + $escaped + '","@")'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" '61.204.47'
Set-TextValue $ws "E2" '  -4.08%  '
Set-TextValue $ws "D3" '2.457.51'
Set-TextValue $ws "E3" '  -6.61%  '
Set-TextValue $ws "D5" '546.72'
Set-TextValue $ws "E5" '  -5.50%  '
Set-TextValue $ws "D6" '146.10'
Set-TextValue $ws "E6" '  -6.86%  '
Set-TextValue $ws "E7" '  +0.01%  '
Set-TextValue $ws "D8" '0.585'
Set-TextValue $ws "E8" '  -7.11%  '
Set-TextValue $ws "D9" '2.457.23'
Set-TextValue $ws "E9" '  -6.63%  '
Set-TextValue $ws "E10" '  -10.24%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws "D11" '0.154'
Set-TextValue $ws "E11" '  -1.72%  '
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws "D12" '5.43'
Set-TextValue $ws "E12" '  -6.62%  '
Set-TextValue $ws "D13" '0.351'
Set-TextValue $ws "E13" '  -8.79%  '
Set-TextValue $ws "D14" '25.98'
Set-TextValue $ws "E14" '  -9.43%  '
Set-TextValue $ws "D15" '2.899.62'
Set-TextValue $ws "E15" '  -6.64%  '
Set-TextValue $ws "E16" '  -9.82%  '
Set-TextValue $ws "D17" '61.135.00'
Set-TextValue $ws "E17" '  -4.08%  '
Set-TextValue $ws "D18" '2.469.24'
Set-TextValue $ws "E18" '  -6.10%  '
Set-TextValue $ws "D19" '11.08'
Set-TextValue $ws "E19" '  -8.95%  '
Set-TextValue $ws "D20" '7.04'
Set-TextValue $ws "E20" '  -8.75%  '
Set-TextValue $ws "D21" '4.15'
Set-TextValue $ws "E21" '  -8.04%  '
Set-TextValue $ws "D22" '318.55'
Set-TextValue $ws "E22" '  -7.31%  '
Set-TextValue $ws "E23" '  -0.03%  '
Set-TextValue $ws "D24" '1.83'
Set-TextValue $ws "E24" '  -3.66%  '
Set-TextValue $ws "D25" '63.83'
Set-TextValue $ws "E25" '  -6.51%  '
Set-TextValue $ws "D26" '2.581.47'
Set-TextValue $ws "E26" '  -6.66%  '
Set-TextValue $ws "D27" '0.0₃0965'
Set-TextValue $ws "E27" '  -14.23%  '
Set-TextValue $ws "D28" '542.27'
Set-TextValue $ws "E28" '  -7.01%  '
Set-TextValue $ws "E29" '  +0.12%  '
Set-TextValue $ws "E30" '  -11.96%  '
Set-TextValue $ws "D31" '8.22'
Set-TextValue $ws "E31" '  -11.08%  '
Set-TextValue $ws "D32" '7.59'
Set-TextValue $ws "E32" '  -8.20%  '
Set-TextValue $ws "E33" '  -8.29%  '
Set-TextValue $ws "E34" '  -7.99%  '
Set-TextValue $ws "E35" '  -9.14%  '
Set-TextValue $ws "D36" '5.83'
Set-TextValue $ws "E36" '  -12.27%  '
Set-TextValue $ws "D37" '1.00'
Set-TextValue $ws "E37" '  +0.11%  '
Set-TextValue $ws "D38" '4.81'
Set-TextValue $ws "E38" '  -12.25%  '
Set-TextValue $ws "D39" '0.377'
Set-TextValue $ws "E39" '  -6.63%  '
Set-TextValue $ws "D40" '18.39'
Set-TextValue $ws "E40" '  -7.04%  '
Set-TextValue $ws "D41" '1.76'
Set-TextValue $ws "E41" '  -8.11%  '
Set-TextValue $ws "D42" '142.08'
Set-TextValue $ws "E42" '  -7.42%  '
Set-TextValue $ws "E43" '  +0.04%  '
Set-TextValue $ws "D44" '40.43'
Set-TextValue $ws "E44" '  -3.86%  '
Set-TextValue $ws "D45" '2.29'
Set-TextValue $ws "E45" '  -10.88%  '
Set-TextValue $ws "D46" '145.60'
Set-TextValue $ws "E46" '  -10.90%  '
Set-TextValue $ws "E47" '  -8.68%  '
Set-TextValue $ws "D48" '21.30'
Set-TextValue $ws "E48" '  -11.74%  '
Set-TextValue $ws "D49" '0.0533'
Set-TextValue $ws "E49" '  -9.18%  '
Set-TextValue $ws "D50" '0.588'
Set-TextValue $ws "E50" '  -7.50%  '
Set-TextValue $ws "D51" '0.0934'
Set-TextValue $ws "E51" '  -6.75%  '
